# Remove the "lemmalist-greek" dependency row from the Acknowledgments sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acknowledgments")

# Row 10 holds the "lemmalist-greek" entry (A10 = "lemmalist-greek").
# Deleting the entire row shifts everything below it up by one, and
# Excel will automatically clean up the now-unused hyperlinks / shared
# strings, and adjust the dimension + sortState ranges.
$ws.Rows.Item(10).Delete()

# Restore the view state (frozen pane top-left cell and active selection)
# to match the saved file.
$ws.Activate()
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D11").Select()
